$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 285, shifting the existing rows
# 285..349 down to 287..351.
$ws.Rows.Item(285).Resize(2).Insert()

# New row 285 data
$ws.Cells.Item(285, 1).Value = 10
$ws.Cells.Item(285, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(285, 3).Value = "La Araucanía"
$ws.Cells.Item(285, 4).Value = 45173
$ws.Cells.Item(285, 5).Value = 9
$ws.Cells.Item(285, 6).Value = 100112013
$ws.Cells.Item(285, 7).Value = "Alcachofa"
$ws.Cells.Item(285, 8).Value = "Española"
$ws.Cells.Item(285, 9).Value = "Primera"
$ws.Cells.Item(285, 10).Value = 200
$ws.Cells.Item(285, 11).Value = 13000
$ws.Cells.Item(285, 12).Value = 13000
$ws.Cells.Item(285, 13).Value = 13000
$ws.Cells.Item(285, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(285, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(285, 16).Value = 433
$ws.Cells.Item(285, 17).Value = 30
$ws.Cells.Item(285, 18).Value = "Hortaliza"

# New row 286 data
$ws.Cells.Item(286, 1).Value = 10
$ws.Cells.Item(286, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(286, 3).Value = "La Araucanía"
$ws.Cells.Item(286, 4).Value = 45173
$ws.Cells.Item(286, 5).Value = 9
$ws.Cells.Item(286, 6).Value = 100112013
$ws.Cells.Item(286, 7).Value = "Alcachofa"
$ws.Cells.Item(286, 8).Value = "Madrigal"
$ws.Cells.Item(286, 9).Value = "Primera"
$ws.Cells.Item(286, 10).Value = 450
$ws.Cells.Item(286, 11).Value = 11000
$ws.Cells.Item(286, 12).Value = 12000
$ws.Cells.Item(286, 13).Value = 11667
$ws.Cells.Item(286, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(286, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(286, 16).Value = 292
$ws.Cells.Item(286, 17).Value = 40
$ws.Cells.Item(286, 18).Value = "Hortaliza"
